# Insert a new data row (row 361) into the worksheet, pushing the existing
# rows 361:480 down to 362:481, and populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 361 (shifts rows 361-480 down to 362-481)
$ws.Rows.Item(361).Insert()

# Populate the newly inserted row 361 with its data
$ws.Cells.Item(361, 1).Value = 10
$ws.Cells.Item(361, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(361, 3).Value = "La Araucanía"
$ws.Cells.Item(361, 4).Value = 44463
$ws.Cells.Item(361, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(361, 5).Value = 9
$ws.Cells.Item(361, 6).Value = 100112006
$ws.Cells.Item(361, 7).Value = "Repollo"
$ws.Cells.Item(361, 8).Value = "Crespo record"
$ws.Cells.Item(361, 9).Value = "Primera"
$ws.Cells.Item(361, 10).Value = 1000
$ws.Cells.Item(361, 11).Value = 900
$ws.Cells.Item(361, 12).Value = 900
$ws.Cells.Item(361, 13).Value = 900
$ws.Cells.Item(361, 14).Value = "$/unidad"
$ws.Cells.Item(361, 15).Value = "Región Metropolitana"
$ws.Cells.Item(361, 16).Value = 900
$ws.Cells.Item(361, 17).Value = 1
$ws.Cells.Item(361, 18).Value = "Hortaliza"
